# "Marked a few strings for translation."
#
# The i18n file list sheet tracks, per source string, whether its
# translation status is "ok" (column B, shared string index 603). A
# block of rows (262-295) had been left without a status; mark them
# "ok" now, matching the style already used for the surrounding B
# column cells (centered, via the existing column style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

for ($r = 262; $r -le 295; $r++) {
    $ws.Cells.Item($r, 2).Value = "ok"
}

# Leave the selection where the author's cursor ended up.
$ws.Range("B295").Select()
